$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.446.77"
$ws.Range("E2").Value = "  -2.47%  "

$ws.Range("D3").Value = "2.979.10"
$ws.Range("E3").Value = "  -4.52%  "

$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$c = $ws.Range("D5")
$c.Value = "'497.09"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.37%  "

$c = $ws.Range("D6")
$c.Value = "'135.77"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.23%  "

$c = $ws.Range("D7")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "2.975.74"
$ws.Range("E8").Value = "  -4.58%  "

$c = $ws.Range("D9")
$c.Value = "'0.427"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "

$c = $ws.Range("D10")
$c.Value = "'7.28"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("E11").Value = "  -3.30%  "

$c = $ws.Range("D12")
$c.Value = "'0.353"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -7.24%  "

$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").Value = "3.493.01"
$ws.Range("E14").Value = "  -4.37%  "

$c = $ws.Range("D15")
$c.Value = "'25.01"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("D16").Value = "56.457.16"
$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").Value = "2.982.97"
$ws.Range("E17").Value = "  -4.34%  "

$c = $ws.Range("D18")
$c.Value = "'0.0000147"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.23%  "

$c = $ws.Range("D19")
$c.Value = "'5.84"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

$c = $ws.Range("D20")
$c.Value = "'12.34"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.92%  "

$c = $ws.Range("D21")
$c.Value = "'7.75"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "

$c = $ws.Range("D22")
$c.Value = "'324.19"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.26%  "

$ws.Range("E23").Value = "  +0.09%  "

$c = $ws.Range("D24")
$c.Value = "'0.463"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -8.37%  "

$c = $ws.Range("D25")
$c.Value = "'61.42"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -10.54%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("D28").Value = "0.0₃0893"
$ws.Range("E28").Value = "  -6.05%  "

$ws.Range("E29").Value = "  +0.02%  "

$c = $ws.Range("D30")
$c.Value = "'6.54"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.29%  "

$c = $ws.Range("D31")
$c.Value = "'6.79"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("E33").Value = "  -6.36%  "

$c = $ws.Range("D34")
$c.Value = "'19.93"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -7.42%  "

$c = $ws.Range("D35")
$c.Value = "'154.38"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.08%  "

$c = $ws.Range("D36")
$c.Value = "'4.51"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -5.00%  "

$c = $ws.Range("D38")
$c.Value = "'5.64"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -8.41%  "

$c = $ws.Range("D39")
$c.Value = "'0.0669"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "

$c = $ws.Range("D40")
$c.Value = "'23.47"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.42%  "

$ws.Range("D41").Value = "3.011.77"
$ws.Range("E41").Value = "  -4.34%  "

$c = $ws.Range("D42")
$c.Value = "'37.20"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -7.88%  "

$c = $ws.Range("D43")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "

$c = $ws.Range("D44")
$c.Value = "'1.01"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -5.84%  "

$ws.Range("E45").Value = "  -0.64%  "

$c = $ws.Range("D46")
$c.Value = "'0.635"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -8.52%  "

$ws.Range("D47").Value = "2.210.37"
$ws.Range("E47").Value = "  -1.84%  "

$c = $ws.Range("D48")
$c.Value = "'3.56"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -8.33%  "

$c = $ws.Range("D49")
$c.Value = "'1.97"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +9.36%  "

$c = $ws.Range("D50")
$c.Value = "'0.0237"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "

$c = $ws.Range("D51")
$c.Value = "'19.30"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.04%  "
